$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 1052
$ws1.Range("F10").Value = 7
$ws1.Range("F13").Value = 13460
$ws1.Range("F17").Value = 5544

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F29").Value = 1052
$ws4.Range("F32").Value = 7
$ws4.Range("F35").Value = 13460
$ws4.Range("F40").Value = 5544
